$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        A = "senior Golang Developer"
        B = "https://www.dice.com/job-detail/9205e164-9988-453d-ab3c-6e306aca3dfa?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "McLean, Virginia"
        D = "Contract, Third Party"
        E = "Depends on Experience"
        F = "NimbusAITech LLC"
    },
    @{
        A = "Golang Developer - San Jose, CA (Preferred) / Remote - 67316"
        B = "https://www.dice.com/job-detail/c25cc8ca-dd38-42ee-ad10-3cde77e958ad?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "Remote or San Jose, California"
        D = "Contract"
        E = '$$40-$43/hr'
        F = "InfiCare Technologies"
    },
    @{
        A = "Go-Lang Developer"
        B = "https://www.dice.com/job-detail/ec6ef353-cf41-4204-a14d-2ab45edb90dd?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "Hybrid in Dallas, Texas"
        D = "Contract"
        E = "Depends on Experience"
        F = "Pyramid Consulting, Inc."
    },
    @{
        A = "Gen AI with Google Cloud Platform"
        B = "https://www.dice.com/job-detail/443ed1e2-c7f1-42dc-8d40-e82b34cca2e5?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang"
        C = "Remote"
        D = "Contract"
        E = "Depends on Experience"
        F = $null
    }
)

$startRow = 159
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($null -ne $row.F) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
}
